# Refresh the cryptocurrency price/volume snapshot (scheduled GitHub Actions
# data pull). Most rows keep their rank but get updated Price (column D) and
# Volume(1h) % change (column E) figures; two pairs of rows (9/10 and 14/15)
# swapped rank order, so their Coin name / Link / Price / Volume cells all
# move together.
#
# Price/volume text is written with NumberFormat "@" (Text) first so Excel's
# automatic type coercion doesn't turn strings like "1.000" or "0.06510"
# into numbers (which would silently drop the formatting the source data
# relies on); the style is then reset to "Normal" so no extra number format
# is left applied to the cell, matching the original (unstyled) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "28.640.75"
    "E2" = "  +3.92%  "
    "D3" = "1.798.39"
    "E3" = "  +0.45%  "
    "D4" = "1.000"
    "E4" = "  -0.03%  "
    "D5" = "313.46"
    "E5" = "  -0.05%  "
    "D7" = "0.5295"
    "E7" = "  -1.49%  "
    "D8" = "0.3781"
    "E8" = "  +0.46%  "
    "B9" = "Dogecoin"
    "C9" = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
    "D9" = "0.07529"
    "E9" = "  +0.25%  "
    "B10" = "OKB"
    "C10" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D10" = "42.65"
    "E10" = "  -0.49%  "
    "D11" = "1.117"
    "E11" = "  +0.42%  "
    "D12" = "1.002"
    "E12" = "  +0.15%  "
    "E13" = "  +1.00%  "
    "B14" = "Polkadot"
    "C14" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D14" = "6.186"
    "E14" = "  +0.30%  "
    "B15" = "Chainlink"
    "C15" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D15" = "7.487"
    "E15" = "  +5.86%  "
    "D16" = "1.796.54"
    "E16" = "  +0.45%  "
    "D17" = "90.30"
    "E17" = "  -0.57%  "
    "E18" = "  -0.41%  "
    "D19" = "0.06468"
    "E19" = "  -0.45%  "
    "D21" = "17.28"
    "E21" = "  +1.86%  "
    "D22" = "5.930"
    "E22" = "  -0.14%  "
    "D23" = "28.636.61"
    "E23" = "  +3.77%  "
    "E24" = "  -0.19%  "
    "D25" = "2.095"
    "E25" = "  +0.46%  "
    "D26" = "161.05"
    "E26" = "  +3.65%  "
    "E27" = "  +0.17%  "
    "D28" = "2.376"
    "D29" = "2.002.26"
    "E29" = "  +0.30%  "
    "D30" = "123.78"
    "E30" = "  +1.58%  "
    "E31" = "  -0.44%  "
    "D32" = "0.1024"
    "E32" = "  -0.71%  "
    "D33" = "5.696"
    "E33" = "  +0.47%  "
    "D34" = "3.683"
    "E34" = "  +2.33%  "
    "D35" = "0.2272"
    "E35" = "  +8.73%  "
    "D36" = "0.06510"
    "E36" = "  +8.11%  "
    "D37" = "8.912"
    "E37" = "  +2.16%  "
    "D38" = "0.02310"
    "E38" = "  +1.04%  "
    "D39" = "5.063"
    "E39" = "  +1.48%  "
    "D40" = "11.45"
    "E40" = "  +0.39%  "
    "D41" = "0.6287"
    "E41" = "  +0.75%  "
    "E42" = "  +5.62%  "
    "D43" = "0.9999"
    "E43" = "  +0.02%  "
    "E44" = "  -1.37%  "
    "D45" = "13.56"
    "E45" = "  +1.68%  "
    "D46" = "0.5912"
    "E46" = "  +0.82%  "
    "E47" = "  +0.79%  "
    "D48" = "126.59"
    "E48" = "  +4.19%  "
    "E49" = "  +3.35%  "
    "D50" = "1.157"
    "E50" = "  +2.12%  "
    "D51" = "0.06928"
    "E51" = "  +2.69%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
